# Record attendance for the 3rd session (column E, week of 2021-02-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance present (1) for rows 3-7, absent (0) for row 8
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 0

# Update the active selection to reflect where the user left off editing
$ws.Range("D14").Select()
